$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.525.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.476.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.44'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.38'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.543'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.505'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.72'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0784'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.111'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.862.16'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.83'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.00'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +8.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.512.73'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.762'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.513.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0936'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.62'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.30'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.31'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.71'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.90'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.85'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.64'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.87'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.21'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.47'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0754'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.46'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -7.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.92'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.105'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.83'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.114'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.11'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.63'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.972.67'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0283'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.14'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.720.23'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.65'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.06'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.30'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.71%  '
